$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.173.06"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.826.26"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'236.20"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").Value = "'0.6040"
$ws.Range("E6").Value = "  -3.97%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.07124"
$ws.Range("E8").Value = "  -5.12%  "

$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").Value = "'24.00"
$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").Value = "1.840.56"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "'4.771"
$ws.Range("E13").Value = "  -4.13%  "

$ws.Range("D14").Value = "'0.6402"
$ws.Range("E14").Value = "  -5.58%  "

$ws.Range("D15").Value = "'0.000009957"
$ws.Range("E15").Value = "  -2.43%  "

$ws.Range("D16").Value = "2.064.99"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").Value = "'79.41"
$ws.Range("E17").Value = "  -3.16%  "

$ws.Range("D18").Value = "'5.986"
$ws.Range("E18").Value = "  -4.05%  "

$ws.Range("D19").Value = "29.214.85"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "'231.53"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "'11.75"
$ws.Range("E22").Value = "  -4.54%  "

$ws.Range("D23").Value = "'7.051"
$ws.Range("E23").Value = "  -4.95%  "

$ws.Range("D24").Value = "'0.9994"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "'155.19"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").Value = "'8.048"
$ws.Range("E26").Value = "  -4.90%  "

$ws.Range("D27").Value = "'0.1284"
$ws.Range("E27").Value = "  -5.09%  "

$ws.Range("D28").Value = "'16.66"
$ws.Range("E28").Value = "  -4.26%  "

$ws.Range("D29").Value = "'0.06841"
$ws.Range("E29").Value = "  +5.44%  "

$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").Value = "'1.450"
$ws.Range("E31").Value = "  -2.36%  "

$ws.Range("D32").Value = "'3.838"
$ws.Range("E32").Value = "  -5.57%  "

$ws.Range("D33").Value = "'3.796"
$ws.Range("E33").Value = "  -6.62%  "

$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("D35").Value = "'1.722"
$ws.Range("E35").Value = "  -6.29%  "

$ws.Range("D36").Value = "'0.6621"
$ws.Range("E36").Value = "  -4.80%  "

$ws.Range("D37").Value = "'2.531"
$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("D38").Value = "1.235.38"
$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("D39").Value = "'2.753"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").Value = "'0.01764"
$ws.Range("E40").Value = "  -4.90%  "

$ws.Range("D41").Value = "'6.588"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("D42").Value = "'0.9335"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").Value = "'0.9997"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "1.987.40"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").Value = "'100.12"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").Value = "'63.40"
$ws.Range("E46").Value = "  -3.29%  "

$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("D48").Value = "'1.638"
$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("D49").Value = "'6.566"
$ws.Range("E49").Value = "  -6.84%  "

# Row 50 and 51 changes (EnergySwap/Cronos reorder)
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05582"
$ws.Range("E50").Value = "  -1.55%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.482"
$ws.Range("E51").Value = "  -5.80%  "

